$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Replace the collaborator name "Keven" with "andrew" in cell A6
$ws.Range("A6").Value = "andrew"

# Update the active selection to A6 to match the new edit location
$ws.Range("A6").Select()
